# Regenerate save_data column G ("K") values for sands_cole.xlsx.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the underlying data pipeline recomputed the
# per-row "K" metric; this script writes the newly computed values into
# column G for rows 2-68 (data rows), leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, in row order starting at row 2 through row 68.
$newK = @(
    1,2,3,0,3,0,1,3,2,1,2,0,0,2,2,1,1,0,2,1,2,0,1,2,2,0,1,3,2,3,0,0,1,0,3,0,
    1,3,1,1,2,1,1,1,0,2,1,0,1,1,1,1,1,3,0,1,2,1,3,4,3,2,0,1,2,1,1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
